$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$data = @(
    @("Danh mục", 0),
    @("Ngày công", 12),
    @("Phụ cấp", 420000),
    @("Lương cơ bản tại CẦN THƠ", 0),
    @("Chiết khấu sale chính tại CẦN THƠ", 0),
    @("Chiết khấu sale phụ tại CẦN THƠ", 0),
    @("Đơn 1 bác sĩ tại CẦN THƠ", 0),
    @("Đơn 2 bác sĩ tại CẦN THƠ", 0),
    @("Công phụ phẫu 1 tại CẦN THƠ", 0),
    @("Công phụ phẫu 2 tại CẦN THƠ", 0),
    @("Lương cơ bản tại LONG XUYÊN", 0),
    @("Chiết khấu sale chính tại LONG XUYÊN", 0),
    @("Chiết khấu sale phụ tại LONG XUYÊN", 0),
    @("Đơn 1 bác sĩ tại LONG XUYÊN", 0),
    @("Đơn 2 bác sĩ tại LONG XUYÊN", 0),
    @("Công phụ phẫu 1 tại LONG XUYÊN", 0),
    @("Công phụ phẫu 2 tại LONG XUYÊN", 0),
    @("Lương cơ bản tại SÓC TRĂNG", 0),
    @("Chiết khấu sale chính tại SÓC TRĂNG", 0),
    @("Chiết khấu sale phụ tại SÓC TRĂNG", 0),
    @("Đơn 1 bác sĩ tại SÓC TRĂNG", 0),
    @("Đơn 2 bác sĩ tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 1 tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 2 tại SÓC TRĂNG", 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
